$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to Text format so numeric-looking values
# ("1.000", "109.04", etc.) are preserved exactly as strings, matching
# the workbook's existing inlineStr storage for that column.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.451.57'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '1.729.71'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '243.76'
$ws.Range("E5").Value = '  +2.70%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = '0.4797'
$ws.Range("E7").Value = '  +3.59%  '
$ws.Range("D8").Value = '0.2667'
$ws.Range("E8").Value = '  +2.88%  '
$ws.Range("D9").Value = '0.06227'
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = '1.729.51'
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").Value = '0.07116'
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '15.72'
$ws.Range("E12").Value = '  +5.31%  '
$ws.Range("D13").Value = '0.6183'
$ws.Range("E13").Value = '  +6.52%  '
$ws.Range("D14").Value = '4.544'
$ws.Range("E14").Value = '  +4.22%  '
$ws.Range("D15").Value = '76.95'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").Value = '26.473.29'
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '0.000006925'
$ws.Range("E19").Value = '  +3.25%  '
$ws.Range("D20").Value = '11.71'
$ws.Range("E20").Value = '  +2.56%  '
$ws.Range("D21").Value = '1.954.46'
$ws.Range("E21").Value = '  +3.97%  '
$ws.Range("D22").Value = '4.561'
$ws.Range("E22").Value = '  +2.14%  '
$ws.Range("D23").Value = '8.893'
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("D24").Value = '5.320'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '136.41'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").Value = '15.35'
$ws.Range("E26").Value = '  +2.21%  '
$ws.Range("D27").Value = '1.791'
$ws.Range("E27").Value = '  +3.64%  '
$ws.Range("D28").Value = '1.405'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").Value = '106.54'
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("D30").Value = '3.984'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = '0.08002'
$ws.Range("E31").Value = '  +4.16%  '
$ws.Range("D32").Value = '3.726'
$ws.Range("E32").Value = '  +3.10%  '
$ws.Range("D33").Value = '0.04550'
$ws.Range("E33").Value = '  +4.49%  '
$ws.Range("D34").Value = '2.614'
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("D35").Value = '0.6402'
$ws.Range("E35").Value = '  +4.62%  '
$ws.Range("D36").Value = '0.9891'
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("D37").Value = '0.9380'
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '1.980'
$ws.Range("E38").Value = '  +5.88%  '
$ws.Range("B39").Value = 'Quant'
$ws.Range("C39").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D39").Value = '107.38'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").Value = '2.408'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").Value = '0.01502'
$ws.Range("E42").Value = '  +3.10%  '
$ws.Range("D43").Value = '5.667'
$ws.Range("E43").Value = '  +11.72%  '
$ws.Range("D44").Value = '0.3906'
$ws.Range("E44").Value = '  +4.68%  '
$ws.Range("D45").Value = '6.942'
$ws.Range("E45").Value = '  +12.81%  '
$ws.Range("D46").Value = '0.1193'
$ws.Range("E46").Value = '  +6.52%  '
$ws.Range("D47").Value = '0.05332'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("D48").Value = '30.79'
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D50").Value = '1.274'
$ws.Range("E50").Value = '  +5.40%  '
$ws.Range("D51").Value = '0.3426'
$ws.Range("E51").Value = '  +2.96%  '

# Clean up: drop the temporary Text number-format override so the
# Price cells go back to carrying no explicit style (as before),
# now that their literal text has been committed.
$ws.Range("D2:D51").Style = "Normal"
